$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = 15
$ws.Range("C10").Value = $ws.Range("C9").Text
$ws.Range("D10").Value = "15. Super Classe para Entidades"
$ws.Range("E10").Value = "4:43 - criação de uma classe abstrata para ser usada como herança, tem a função de eliminar a necessidade de declarar/setar ids nas entidades. Para seu uso, ao implementar a entidade é necessário colocar a instrução ""extends"" herdando assim a classe abstrata passando como parametro um tipo LONG"

$ws.Range("E10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 60

$ws.Range("C9").Select()
